$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the self-evaluation ("auto-evaluation") scores in column D for each
# criterion row of the evaluation grid.
$values = @{
    3  = 0
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 1
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 0.5
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    22 = 1
    23 = 1
    24 = 0.5
    25 = 1
    29 = 0.5
    30 = 0
    31 = 1
    32 = 0
    33 = 0
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    41 = 1
    42 = 0
    43 = 0
    44 = 0
    45 = 0
    46 = 0
    47 = 0
    48 = 0
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}

# Move the active selection to D49 (scrolls the view back to the top as a
# side effect, matching the saved workbook state).
$ws.Range("D49").Select() | Out-Null

# Re-create the (hidden) workbook-level defined name that was present in the
# saved workbook.
$definedName = $wb.Names.Add("LOCAL_MYSQL_DATE_FORMAT", "=REPT(LOCAL_YEAR_FORMAT,4)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_MONTH_FORMAT,2)&LOCAL_DATE_SEPARATOR&REPT(LOCAL_DAY_FORMAT,2)&"" ""&REPT(LOCAL_HOUR_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_MINUTE_FORMAT,2)&LOCAL_TIME_SEPARATOR&REPT(LOCAL_SECOND_FORMAT,2)")
$definedName.Visible = $false | Out-Null
